# Spring 2026 (1261) term update for the LS (Liberal Studies) department
# requisites sheet:
#   - Remove retired/duplicate course rows: LS991, LS898, LS990
#   - Re-sort the remaining course rows alphabetically/numerically by
#     Course_Code
#   - Refresh the selection/used-range to the new, smaller extent

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the rows that were removed for this term. Work bottom-up so the
# row numbers of the rows still to be removed don't shift underneath us.
# Row 16 = LS990, Row 15 = LS898, Row 7 = LS991 (original layout).
$ws.Rows(16).Delete()
$ws.Rows(15).Delete()
$ws.Rows(7).Delete()

# Re-sort the remaining course rows (now A2:I15) by Course_Code.
$dataRange = $ws.Range("A2:I15")
$dataRange.Sort($ws.Range("A2:A15"))

# Match the refreshed selection left behind by the edit.
$ws.Range("A2:I15").Select()
